$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.920.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "'1.899.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'0.7967"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.31%  "

$ws.Range("D6").Value = "'244.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.3180"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.95%  "

$ws.Range("D9").Value = "'25.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.17%  "

$ws.Range("D10").Value = "'0.07208"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").Value = "'0.08117"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'5.642"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.10%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7734"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").Value = "'1.886.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "'92.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "'6.222"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.84%  "

$ws.Range("D17").Value = "'29.901.65"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'245.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").Value = "'0.000007788"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'8.251"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +18.13%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'2.144.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  -4.14%  "

$ws.Range("D26").Value = "'9.494"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.47%  "

$ws.Range("D27").Value = "'164.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").Value = "'18.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("D29").Value = "'2.079"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.49%  "

$ws.Range("E30").Value = "  +3.48%  "

$ws.Range("D31").Value = "'1.552"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("D32").Value = "'4.509"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.88%  "

$ws.Range("D33").Value = "'0.05658"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.57%  "

$ws.Range("D34").Value = "'4.092"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").Value = "'1.292"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

$ws.Range("D36").Value = "'0.7474"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.10%  "

$ws.Range("D37").Value = "'1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("E38").Value = "  -2.89%  "

$ws.Range("E39").Value = "  +0.96%  "

$ws.Range("D40").Value = "'2.791"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").Value = "'1.171.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.12%  "

$ws.Range("D42").Value = "'74.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.06%  "

$ws.Range("D43").Value = "'0.4442"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'5.967"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("D45").Value = "'0.8571"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("D46").Value = "'104.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "'10.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.40%  "

$ws.Range("D49").Value = "'1.889"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "'7.513"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("D51").Value = "'2.965"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.28%  "
